# "Generate Report for handback"
#
# The handback-status report keeps one row per source file per locale
# sheet. Re-running the report for the first source file
# (3120aef5-5742-44c4-bfc8-f48d3381e7be...) picked up a fresh handback
# pass, so its "Correspond Handoff Datetime" (column D) and "Correspond
# Handback DateTime" (column G) timestamps advance on row 2 of each
# locale sheet. The second source file's row (a1eefbe0-7ce3-406d...,
# row 3) was not part of this pass, so its timestamps stay the same.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 = 3120aef5-5742-44c4-bfc8-f48d3381e7be... file
$wsZhCn.Range("D2").Value = "2016-01-17 03:20:01"
$wsZhCn.Range("G2").Value = "2016-01-17 03:20:44"

# zh-cn sheet: row 3 = a1eefbe0-7ce3-406d-afc9-ca13d215af41... file (unchanged)
$wsZhCn.Range("D3").Value = "2016-01-17 03:18:19"
$wsZhCn.Range("G3").Value = "2016-01-17 03:19:00"

# de-de sheet: row 2 = 3120aef5-5742-44c4-bfc8-f48d3381e7be... file
$wsDeDe.Range("D2").Value = "2016-01-17 03:20:11"
$wsDeDe.Range("G2").Value = "2016-01-17 03:21:00"

# de-de sheet: row 3 = a1eefbe0-7ce3-406d-afc9-ca13d215af41... file (unchanged)
$wsDeDe.Range("D3").Value = "2016-01-17 03:18:30"
$wsDeDe.Range("G3").Value = "2016-01-17 03:19:17"
